$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "valid as of" date used in the regression test sheet
$ws.Range("O2").Value = "15/04/2024"

# Update the monthly period references from 202305 to 202405
$ws.Range("Q2").Value = "202405"
$ws.Range("T2").Value = 202405

# Move the active selection (matches the recorded view state in the diff)
$ws.Application.ActiveWindow.ScrollColumn = 12
$ws.Range("T3").Select()
